# Updated cryptos list: refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.830.98"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "1.641.17"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("E4").Value = "  -0.66%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  +1.55%  "
$ws.Range("E7").Value = "  -0.69%  "
$ws.Range("E8").Value = "  +1.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0620"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0844"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").Value = "1.871.29"
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("D13").Value = "1.630.51"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("E15").Value = "  +1.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.60%  "
$ws.Range("D17").Value = "26.851.21"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").Value = "0.0₃0730"
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.64%  "
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.62%  "
$ws.Range("E22").Value = "  +1.37%  "
$ws.Range("E23").Value = "  +3.38%  "
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.30%  "
$ws.Range("E26").Value = "  -0.80%  "
$ws.Range("E27").Value = "  +4.86%  "
$ws.Range("E28").Value = "  +1.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0510"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.95%  "
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("E33").Value = "  +2.06%  "
$ws.Range("E34").Value = "  +2.61%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "1.235.24"
$ws.Range("E36").Value = "  -2.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0174"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("E38").Value = "  +3.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.831"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.98%  "
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.807"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.32%  "
$ws.Range("D43").Value = "1.783.86"
$ws.Range("E43").Value = "  +0.66%  "
$ws.Range("E44").Value = "  -2.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "60.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.53%  "
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("E48").Value = "  +11.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0513"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0971"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.02%  "
